# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: AD=30, AE=31, AF=32
$colWins   = 30
$colLosses = 31
$colTies   = 32

# --- Header row (row 1), styled like the other headers (bold, centered, bordered) ---
$headerCells = @(
    @{ Col = $colWins;   Text = "Wins" },
    @{ Col = $colLosses; Text = "Losses" },
    @{ Col = $colTies;   Text = "Ties" }
)

foreach ($h in $headerCells) {
    $cell = $ws.Cells.Item(1, $h.Col)
    $cell.Value = $h.Text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
}

# --- Data rows (2 through 49): every team entry gets the same season record ---
$wins = 83
$losses = 79
$ties = 0

for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, $colWins).Value = $wins
    $ws.Cells.Item($r, $colLosses).Value = $losses
    $ws.Cells.Item($r, $colTies).Value = $ties
}
